$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the date a record was last changed.
# Every data row (2 through 277) has its "changed" date advanced by one day,
# from serial date 46074 (2026-02-21) to 46075 (2026-02-22).
for ($r = 2; $r -le 277; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = 46075
}
